$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row values (B2:E2)
$ws.Range("B2").Value = 5.3266291156268153
$ws.Range("C2").Value = 6.905573518806766
$ws.Range("D2").Value = 3.8708698498221517
$ws.Range("E2").Value = 4.2224481990286042

# Update STR row values (B3:E3)
$ws.Range("B3").Value = 8.0717059160327356
$ws.Range("C3").Value = 13.749156295846296
$ws.Range("D3").Value = 10.79845978970932
$ws.Range("E3").Value = 3.1364440778250153

# Update selection to match new reduced range
$ws.Range("B1:E3").Select()
